{"js": "const body = context.document.body;\nconst results = body.search(\"Manuel Dias\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\nif (results.items.length > 0) {\n  results.items[0].insertText(\"Manuel In\u00e1cio Veladas Dias\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"Manuel Dias\", $false, $false, $false, $false, $false, $true, 1, $false, \"Manuel In\u00e1cio Veladas Dias\", 2)\n"}
